$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 'General and Preventive Care, Mental Health Services'
$ws.Range("Q3").Value = 'Diagnostic and Imaging Services, Rehabilitation and Therapy, Treatment and Procedures, Pediatric Care'
$ws.Range("Q4").Value = 'Rehabilitation and Therapy, Emergency and Critical Care, Treatment and Procedures'
$ws.Range("Q5").Value = 'Women''s Health, Treatment and Procedures, Neonatal Care, Diagnostic and Imaging Services, Emergency and Critical Care, Oncology'
$ws.Range("Q6").Value = 'Oncology, Palliative and Supportive Care'
$ws.Range("Q9").Value = 'Women''s Health, Orthopedics, Treatment and Procedures, Pediatric Care, Diagnostic and Imaging Services, Emergency and Critical Care, General and Preventive Care'
$ws.Range("Q18").Value = 'Ophthalmology, Orthopedics, Rehabilitation and Therapy, Treatment and Procedures, Mental Health Services, Dermatology, ENT (Ear, Nose, Throat), Diagnostic and Imaging Services, Emergency and Critical Care, Infectious Diseases, General and Preventive Care, Patient Care'
$ws.Range("Q25").Value = 'Emergency and Critical Care, Mental Health Services, Patient Care'
$ws.Range("Q36").Value = 'Ophthalmology, Women''s Health, Orthopedics, Rehabilitation and Therapy, Neonatal Care, Emergency and Critical Care, Infectious Diseases, Patient Care'
$ws.Range("Q38").Value = 'Women''s Health, Treatment and Procedures, Infectious Diseases, General and Preventive Care, Patient Care'
$ws.Range("Q40").Value = 'Treatment and Procedures, Emergency and Critical Care, General and Preventive Care'
$ws.Range("Q41").Value = 'Infectious Diseases, Treatment and Procedures, Patient Care, Women''s Health'
$ws.Range("Q42").Value = 'General and Preventive Care, Women''s Health'
$ws.Range("Q43").Value = 'Emergency and Critical Care, General and Preventive Care'
$ws.Range("Q44").Value = 'Ophthalmology, Orthopedics, Diagnostic and Imaging Services, Infectious Diseases, Patient Care'
$ws.Range("Q48").Value = 'Rehabilitation and Therapy, Patient Care, Treatment and Procedures, Women''s Health'
